$wb = $excel.ActiveWorkbook

# "Overview" sheet: the 889e6598 file's status moves from
# "Handed back: in sync with en-US" to "Ready for handoff".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# "zh-cn" sheet: same status change, plus a refreshed handoff datetime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-01-18 05:44:19"

# "de-de" sheet: same status change, plus a refreshed handoff datetime.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-01-18 05:44:32"
